$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.946.54"
$ws.Range("E2").Value = "  -4.02%  "
$ws.Range("D3").Value = "'2.374.63"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'292.50"
$ws.Range("E5").Value = "  -3.44%  "
$ws.Range("D6").Value = "'93.24"
$ws.Range("E6").Value = "  -8.07%  "
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").Value = "'33.87"
$ws.Range("E10").Value = "  -5.49%  "
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("E12").Value = "  -3.93%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "'2.739.37"
$ws.Range("E14").Value = "  +4.46%  "
$ws.Range("D15").Value = "'2.375.85"
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").Value = "'13.85"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "'0.818"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "'44.973.20"
$ws.Range("E18").Value = "  -3.94%  "
$ws.Range("D19").Value = "'12.36"
$ws.Range("E19").Value = "  -4.97%  "
$ws.Range("D20").Value = "'0.0₃0926"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "'6.06"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").Value = "'66.12"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").Value = "'236.90"
$ws.Range("E23").Value = "  -5.18%  "
$ws.Range("D24").Value = "'2.74"
$ws.Range("E24").Value = "  -3.99%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").Value = "'37.04"
$ws.Range("E28").Value = "  -13.21%  "
$ws.Range("D29").Value = "'9.48"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").Value = "'3.78"
$ws.Range("E30").Value = "  +14.06%  "
$ws.Range("D31").Value = "'20.78"
$ws.Range("E31").Value = "  +4.41%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'146.70"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.68"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("D34").Value = "'5.36"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").Value = "'0.0753"
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").Value = "'1.92"
$ws.Range("E37").Value = "  +10.53%  "
$ws.Range("D38").Value = "'0.113"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "'14.46"
$ws.Range("E39").Value = "  -10.68%  "
$ws.Range("D40").Value = "'3.68"
$ws.Range("E40").Value = "  -6.30%  "
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").Value = "'1.956.13"
$ws.Range("E42").Value = "  +8.02%  "
$ws.Range("D43").Value = "'3.13"
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'87.91"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").Value = "'1.70"
$ws.Range("E46").Value = "  -14.55%  "
$ws.Range("D47").Value = "'8.36"
$ws.Range("E47").Value = "  +6.52%  "
$ws.Range("E48").Value = "  +15.59%  "
$ws.Range("D49").Value = "'98.60"
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("D50").Value = "'2.609.73"
$ws.Range("E50").Value = "  +4.43%  "
$ws.Range("D51").Value = "'0.180"
$ws.Range("E51").Value = "  -5.41%  "
